$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.760.55"

$ws.Range("D3").Value = "2.440.18"
$ws.Range("E3").Value = "  -1.42%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'560.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").Value = "'163.34"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -1.01%  "

$ws.Range("D9").Value = "'0.170"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.26%  "

$ws.Range("E10").Value = "  -2.24%  "

$ws.Range("E11").Value = "  -0.66%  "

$ws.Range("D12").Value = "'4.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.19%  "

$ws.Range("E13").Value = "  +4.10%  "

$ws.Range("D14").Value = "68.659.06"
$ws.Range("E14").Value = "  -0.57%  "

$ws.Range("D15").Value = "2.887.04"
$ws.Range("E15").Value = "  -0.64%  "

$ws.Range("D16").Value = "'23.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.60%  "

$ws.Range("D17").Value = "2.442.15"
$ws.Range("E17").Value = "  -2.56%  "

$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("D19").Value = "'338.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "'7.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.51%  "

$ws.Range("D21").Value = "'3.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("E22").Value = "  +2.15%  "

$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").Value = "'65.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.87%  "

$ws.Range("D25").Value = "'3.78"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.23%  "

$ws.Range("D26").Value = "2.566.97"
$ws.Range("E26").Value = "  -1.46%  "

$ws.Range("D27").Value = "'8.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.78%  "

$ws.Range("D28").Value = "'0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("D29").Value = "0.0₃0823"
$ws.Range("E29").Value = "  -1.00%  "

$ws.Range("D30").Value = "'7.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.05%  "

$ws.Range("E31").Value = "  +3.55%  "

$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").Value = "'432.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.15%  "

$ws.Range("E34").Value = "  -2.24%  "

$ws.Range("D35").Value = "'158.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.34%  "

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("D38").Value = "'17.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("E39").Value = "  -0.78%  "

$ws.Range("D40").Value = "'0.301"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.41%  "

$ws.Range("E41").Value = "  +1.66%  "

$ws.Range("E42").Value = "  -1.70%  "

$ws.Range("E43").Value = "  -0.37%  "

$ws.Range("D44").Value = "'2.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("E45").Value = "  -1.78%  "

$ws.Range("D46").Value = "'130.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "

$ws.Range("D47").Value = "'0.0720"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.35%  "

$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("E49").Value = "  -1.78%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.0923"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.33%  "

$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").Value = "'1.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.65%  "
